$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Main scene inventory values (row 4) — Summer/Winter/WindChance/SnowChance updates
$ws.Range("F4").Value = 2
$ws.Range("H4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0

# Give the Spring/Summer/Fall/Winter columns explicit widths so the data is readable
$ws.Columns("E").ColumnWidth = 14.035714285714286
$ws.Columns("F").ColumnWidth = 11.535714285714286
$ws.Columns("G").ColumnWidth = 13.660714285714286
$ws.Columns("H").ColumnWidth = 11.660714285714286

# Move the active selection (path fix) to N10
$ws.Range("N10").Select() | Out-Null
